$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear previous contents/formatting so the sheet can be rebuilt cleanly ---
$ws.UsedRange.Clear()

# --- Column widths (approximation of original pixel-based "best fit" widths) ---
$ws.Columns.Item(1).ColumnWidth = 10.451822916666666
$ws.Columns.Item(2).ColumnWidth = 4.451822916666667
$ws.Columns.Item(3).ColumnWidth = 4.877604166666667
$ws.Columns.Item(4).ColumnWidth = 5.307291666666667
$ws.Columns.Item(5).ColumnWidth = 6.166666666666667
$ws.Columns.Item(6).ColumnWidth = 6.592447916666667
$ws.Columns.Item(7).ColumnWidth = 6.592447916666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 6.736979166666667
$ws.Columns.Item(10).ColumnWidth = 7.307291666666667
$ws.Columns.Item(11).ColumnWidth = 4.736979166666667
$ws.Columns.Item(12).ColumnWidth = 5.307291666666667
$ws.Columns.Item(14).ColumnWidth = 15.592447916666666
$ws.Columns.Item(15).ColumnWidth = 19.166666666666668

# --- Header row (row 2) text ---
$ws.Range("A2").Value = 'modelo'
$ws.Range("B2").Value = 'bos_token'
$ws.Range("C2").Value = 'eos_token'
$ws.Range("D2").Value = 'sep_token'
$ws.Range("E2").Value = 'cls_token'
$ws.Range("F2").Value = 'pad_token'
$ws.Range("G2").Value = 'padding_side'
$ws.Range("H2").Value = 'mask_token'
$ws.Range("I2").Value = 'unk_token'
$ws.Range("J2").Value = 'subtokens'
$ws.Range("K2").Value = 'subtoken'
$ws.Range("L2").Value = 'lower_case'
$ws.Range("M2").Value = 'algoritmo'
$ws.Range("N2").Value = 'sentença simples'
$ws.Range("O2").Value = 'par de sentenças'
$ws.Range("P2").Value = 'Exemplo'

# --- Header styling (bold, size 12) - done early so this style claims cellXfs index 2 ---
$headerRange = $ws.Range("A2:P2")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$ws.Rows.Item(2).RowHeight = 15.75

# --- Row 3: BERT ---
$ws.Range("A3").Value = 'BERT'
$ws.Range("B3").Value = '[CLS]'
$ws.Range("C3").Value = '[SEP]'
$ws.Range("D3").Value = '[SEP]'
$ws.Range("E3").Value = '[CLS]'
$ws.Range("F3").Value = '[PAD]'
$ws.Range("G3").Value = 'direita'
$ws.Range("H3").Value = '[MASK]'
$ws.Range("I3").Value = '[UNK]'
$ws.Range("K3").Value = '##'
$ws.Range("M3").Value = 'Wordpiece'
$ws.Range("N3").Value = '[CLS] X [SEP]'
$ws.Range("O3").Value = '[CLS] A [SEP] B [SEP]'
$ws.Range("P3").Value = '[''[CLS]'',''Ado'', ''##ro'', ''sor'', ''##vete'', ''de'', ''mang'', ''##a'', ''.'', ''[SEP]'']'

# --- Row 4: Albert ---
$ws.Range("A4").Value = 'Albert'
$ws.Range("B4").Value = '[CLS]'
$ws.Range("C4").Value = '[SEP]'
$ws.Range("D4").Value = '[SEP]'
$ws.Range("E4").Value = '[CLS]'
$ws.Range("F4").Value = '<pad>'
$ws.Range("G4").Value = 'direita'
$ws.Range("H4").Value = '[MASK]'
$ws.Range("I4").Value = '<unk>'
$ws.Range("K4").Value = '▁'
$ws.Range("M4").Value = 'SentencePiece'
$ws.Range("N4").Value = '[CLS] X [SEP]'
$ws.Range("O4").Value = '[CLS] A [SEP] B [SEP]'
$ws.Range("P4").Value = '[''[CLS]'',''▁a'', ''doro'', ''▁sor'', ''ve'', ''te'', ''▁de'', ''▁manga'', ''.'',''[SEP]'']'

# --- Row 5: Distilbert ---
$ws.Range("A5").Value = 'Distilbert'
$ws.Range("B5").Value = '[CLS]'
$ws.Range("C5").Value = '[SEP]'
$ws.Range("D5").Value = '[SEP]'
$ws.Range("E5").Value = '[CLS]'
$ws.Range("F5").Value = '[PAD]'
$ws.Range("G5").Value = 'direita'
$ws.Range("H5").Value = '[MASK]'
$ws.Range("I5").Value = '[UNK]'
$ws.Range("M5").Value = 'Wordpice'
$ws.Range("N5").Value = '<s> X </s>'
$ws.Range("O5").Value = '<s> A </s></s> B </s>'
$ws.Range("P5").Value = '[''[CLS]'',''I'', ''play'', ''bass'', ''in'', ''a'', ''jazz'', ''band'', ''.'',''[SEP]'']'

# --- Row 6: RoBERTa ---
$ws.Range("A6").Value = 'RoBERTa'
$ws.Range("B6").Value = '<s>'
$ws.Range("C6").Value = '</s>'
$ws.Range("D6").Value = '</s>'
$ws.Range("E6").Value = '<s>'
$ws.Range("F6").Value = '<pad>'
$ws.Range("G6").Value = 'direita'
$ws.Range("H6").Value = '<mask>'
$ws.Range("I6").Value = 'Â'
$ws.Range("K6").Value = 'Ġ'
$ws.Range("M6").Value = 'byte-pair-encoding (BPE) [ Sennrich et al. ] '
$ws.Range("N6").Value = '[CLS] X [SEP]'
$ws.Range("O6").Value = '[CLS] A [SEP] B [SEP]'
$ws.Range("P6").Value = '[''<s>'',''Su'', ''je'', ''i'', ''Ġa'', ''Ġmanga'', ''Ġda'', ''Ġcam'', ''isa'', ''.'',''</s>]'

# --- Row 7: XLNet ---
$ws.Range("A7").Value = 'XLNet'
$ws.Range("B7").Value = '<s>'
$ws.Range("C7").Value = '</s>'
$ws.Range("D7").Value = '<sep>'
$ws.Range("E7").Value = '<cls>'
$ws.Range("F7").Value = '<pad>'
$ws.Range("G7").Value = 'esquerda'
$ws.Range("H7").Value = '<mask>'
$ws.Range("I7").Value = '<unk>'
$ws.Range("K7").Value = '▁'
$ws.Range("M7").Value = 'SentencePiece'
$ws.Range("N7").Value = 'X <sep> <cls>'
$ws.Range("O7").Value = 'A <sep> B <sep> <cls>'
$ws.Range("P7").Value = '[''▁a'', ''doro'', ''▁sor'', ''ve'', ''te'', ''▁de'', ''▁manga'', ''.'',''<sep>'',''<cls>'']'

# --- Bold model-name column (A3:A7), matching the pre-existing bold style used for column A ---
$ws.Range("A3:A7").Font.Bold = $true

# --- Force the boolean-looking "true"/"false" values to be stored as literal text        ---
# --- (leading apostrophe forces text entry, then resetting the style drops the visual    ---
# --- quote-prefix indicator so the cells end up unstyled, matching plain text cells)      ---
$ws.Range("J3").Value = "'true"
$ws.Range("J3").Style = "Normal"
$ws.Range("J4").Value = "'true"
$ws.Range("J4").Style = "Normal"
$ws.Range("L4").Value = "'true"
$ws.Range("L4").Style = "Normal"
$ws.Range("J5").Value = "'false"
$ws.Range("J5").Style = "Normal"
$ws.Range("J6").Value = "'true"
$ws.Range("J6").Style = "Normal"
$ws.Range("J7").Value = "'true"
$ws.Range("J7").Style = "Normal"

# --- Selection ---
$ws.Range("G8").Select()

Write-Host "Done"
